$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 1.3
$ws.Range("L2").Value = 1.17
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 2.88
$ws.Range("O2").Value = 1.08
$ws.Range("P2").Value = 2.88
$ws.Range("Q2").Value = 1.08
$ws.Range("R2").Value = 1.87
$ws.Range("S2").Value = 1.08
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 4.4
$ws.Range("W2").Value = 1.02
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
$ws.Range("F3").Value = 2.54
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 2.96
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 3.05
$ws.Range("K3").Value = 3.5
$ws.Range("P3").Value = 1.57
$ws.Range("G4").Value = 1.45
$ws.Range("H4").Value = 8
$ws.Range("J4").Value = 4.3
$ws.Range("L4").Value = 1.25
$ws.Range("Q4").Value = 1.48
$ws.Range("R4").Value = 1.66
$ws.Range("W4").Value = 3.4
$ws.Range("AN4").Value = 4.7
$ws.Range("I5").Value = 1.86
$ws.Range("F6").Value = 1.81
$ws.Range("G6").Value = 1.96
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 3.7
$ws.Range("K6").Value = 4.3
$ws.Range("G7").Value = 5.6
$ws.Range("J7").Value = 3.65
$ws.Range("P7").Value = 1.78
$ws.Range("Q7").Value = 1.01
$ws.Range("F8").Value = 1.04
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.01
$ws.Range("F10").Value = 9
$ws.Range("I10").Value = 1.4
$ws.Range("J10").Value = 5.6
$ws.Range("K10").Value = 5.7
$ws.Range("P10").Value = 2.02
$ws.Range("Q10").Value = 1.69
$ws.Range("F11").Value = 3.2
$ws.Range("I11").Value = 1.94
$ws.Range("J11").Value = 4
$ws.Range("Q11").Value = 1.38
$ws.Range("F12").Value = 1.09
$ws.Range("G12").Value = 1000
$ws.Range("H12").Value = 1.04
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 1.01
$ws.Range("K12").Value = 7.2
$ws.Range("P12").Value = 1.15
$ws.Range("Q12").Value = 1.01
$ws.Range("H13").Value = 1.84
$ws.Range("I13").Value = 2
$ws.Range("P13").Value = 2.74
$ws.Range("Q13").Value = 1.47
$ws.Range("I14").Value = 3.45
$ws.Range("J14").Value = 2.9
$ws.Range("P14").Value = 2.58
$ws.Range("Q14").Value = 1.47
$ws.Range("F15").Value = 1.04
$ws.Range("H15").Value = 1.04
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 1.01
$ws.Range("F17").Value = 1.04
$ws.Range("H17").Value = 1.04
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1.01
$ws.Range("K17").Value = 60
$ws.Range("P17").Value = 1.24
$ws.Range("Q17").Value = 1.01
$ws.Range("F18").Value = 2.18
$ws.Range("H18").Value = 1.53
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 1.53
$ws.Range("K18").Value = 1000
$ws.Range("Q18").Value = 1.01
$ws.Range("G19").Value = 4.8
$ws.Range("H19").Value = 1.97
$ws.Range("N20").Value = 3.15
$ws.Range("P20").Value = 1.72
$ws.Range("T20").Value = 2.2
$ws.Range("Z20").Value = 44
$ws.Range("AD20").Value = 23
$ws.Range("AE20").Value = 1000
$ws.Range("AF20").Value = 8.800000000000001
$ws.Range("F21").Value = 2.6
$ws.Range("G21").Value = 2.64
$ws.Range("I21").Value = 3.25
$ws.Range("Q21").Value = 2.36
$ws.Range("U21").Value = 1.97
$ws.Range("F22").Value = 2.1
$ws.Range("H22").Value = 4.2
$ws.Range("I22").Value = 4.8
$ws.Range("J22").Value = 2.94
$ws.Range("K22").Value = 3.2
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = 1.27
$ws.Range("I23").Value = 1.31
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = 7
$ws.Range("P23").Value = 2.02
$ws.Range("Q23").Value = 1.01
$ws.Range("G24").Value = 1000
$ws.Range("H24").Value = 1.36
$ws.Range("P24").Value = 1.73
$ws.Range("Q24").Value = 1.97
$ws.Range("F26").Value = 3.3
$ws.Range("G26").Value = 4.2
$ws.Range("H26").Value = 2.36
$ws.Range("I26").Value = 2.7
$ws.Range("F27").Value = 3.45
$ws.Range("F30").Value = 1.04
$ws.Range("G30").Value = 1000
$ws.Range("H30").Value = 1.04
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 1.01
$ws.Range("K30").Value = 1000
$ws.Range("P30").Value = 1.24
$ws.Range("Q30").Value = 1.01
